# Push clean data of population_age_sex_race.csv and state_name_clean at resources
# Adds a new "87+" age bracket / "lost" generation row (row 12) to the
# AGE / GENERATION table on Sheet1, mirroring the existing A1..A5 / G1..G5
# rows already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "A6"
$ws.Range("B12").Value = "87+"
$ws.Range("C12").Value = "G6"
$ws.Range("D12").Value = "A6"
$ws.Range("E12").Value = "lost"

# Leave the selection on the newly-entered cell, matching the author's
# final cursor position.
$ws.Range("E12").Select() | Out-Null
